# Updated cryptos list (data refresh) - applies the latest price/volume snapshot
# to the Sheet1 table (coin prices in column D, 1h volume % in column E, plus a
# ranking swap between Avalanche/ShibaInu rows 13-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) cells: force text storage so number-like strings
# (e.g. "1.00", "0.999", "11.30") are not reinterpreted as numeric values,
# then restore the default "Normal" style so no stray formatting is left behind.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.645.51'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.441.83'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.85'
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.534'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.112'
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.27'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.353'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '26.94'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000184'
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.537.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.437.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '325.02'
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.87'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '566.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0₃0990'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.556.43'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.148'
$ws.Range("D32").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.91'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.383'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.52'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.78'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '150.61'
$ws.Range("D40").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '149.62'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0538'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.51'
$ws.Range("D47").Style = "Normal"

# --- Other changed cells (Coin name, Link, Volume(1h)) -- plain text assignment.
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("E3").Value = '  +1.15%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("E6").Value = '  +2.22%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  +0.81%  '
$ws.Range("E9").Value = '  +2.40%  '
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("E13").Value = '  +5.23%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("E14").Value = '  +5.29%  '
$ws.Range("E16").Value = '  +0.82%  '
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("E18").Value = '  +0.17%  '
$ws.Range("E19").Value = '  +2.02%  '
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  +2.65%  '
$ws.Range("E24").Value = '  +3.53%  '
$ws.Range("E25").Value = '  -1.55%  '
$ws.Range("E26").Value = '  -1.76%  '
$ws.Range("E27").Value = '  +4.61%  '
$ws.Range("E28").Value = '  +1.00%  '
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("E30").Value = '  +2.20%  '
$ws.Range("E31").Value = '  +2.16%  '
$ws.Range("E32").Value = '  -0.25%  '
$ws.Range("E33").Value = '  +1.00%  '
$ws.Range("E34").Value = '  +1.34%  '
$ws.Range("E35").Value = '  +4.68%  '
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  -0.21%  '
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("E40").Value = '  -0.91%  '
$ws.Range("E41").Value = '  +0.76%  '
$ws.Range("E42").Value = '  +0.24%  '
$ws.Range("E43").Value = '  +5.58%  '
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("E45").Value = '  +1.48%  '
$ws.Range("E46").Value = '  +1.04%  '
$ws.Range("E47").Value = '  +2.58%  '
$ws.Range("E48").Value = '  +1.63%  '
$ws.Range("E49").Value = '  +1.65%  '
$ws.Range("E50").Value = '  +1.92%  '
$ws.Range("E51").Value = '  +0.59%  '
